$d = $word.ActiveDocument

# Update the date/weekday heading paragraph
$d.Paragraphs.Item(1).Range.Find.Execute("2024-01-07 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-01-08 Monday", 2) | Out-Null

# Update the division problems in the table, cell by cell to avoid cross-matching
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Find.Execute("49÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷9=", 2) | Out-Null
$t.Cell(1, 2).Range.Find.Execute("80÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷2=", 2) | Out-Null
$t.Cell(1, 3).Range.Find.Execute("63÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷3=", 2) | Out-Null
$t.Cell(1, 4).Range.Find.Execute("34÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=", 2) | Out-Null
$t.Cell(1, 5).Range.Find.Execute("70÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷2=", 2) | Out-Null
$t.Cell(5, 1).Range.Find.Execute("39÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=", 2) | Out-Null
$t.Cell(5, 2).Range.Find.Execute("24÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷3=", 2) | Out-Null
$t.Cell(5, 3).Range.Find.Execute("18÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷9=", 2) | Out-Null
$t.Cell(5, 4).Range.Find.Execute("42÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷5=", 2) | Out-Null
$t.Cell(5, 5).Range.Find.Execute("77÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷9=", 2) | Out-Null
$t.Cell(9, 1).Range.Find.Execute("17÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=", 2) | Out-Null
$t.Cell(9, 2).Range.Find.Execute("17÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷4=", 2) | Out-Null
$t.Cell(9, 3).Range.Find.Execute("85÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "87÷6=", 2) | Out-Null
$t.Cell(9, 4).Range.Find.Execute("77÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "90÷5=", 2) | Out-Null
$t.Cell(9, 5).Range.Find.Execute("14÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "53÷8=", 2) | Out-Null
$t.Cell(13, 1).Range.Find.Execute("97÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "85÷4=", 2) | Out-Null
$t.Cell(13, 2).Range.Find.Execute("67÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷7=", 2) | Out-Null
$t.Cell(13, 3).Range.Find.Execute("26÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷8=", 2) | Out-Null
$t.Cell(13, 4).Range.Find.Execute("56÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷4=", 2) | Out-Null
$t.Cell(13, 5).Range.Find.Execute("11÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷7=", 2) | Out-Null
$t.Cell(17, 1).Range.Find.Execute("55÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=", 2) | Out-Null
$t.Cell(17, 2).Range.Find.Execute("11÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷7=", 2) | Out-Null
$t.Cell(17, 3).Range.Find.Execute("20÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷3=", 2) | Out-Null
$t.Cell(17, 4).Range.Find.Execute("67÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "94÷9=", 2) | Out-Null
$t.Cell(17, 5).Range.Find.Execute("79÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷3=", 2) | Out-Null
